$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range / last row with data
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Swap the contents of column C (codeforiati:group-code) and
# column D (codeforiati:group-name) for every row, including the header.
for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $dVal
    $ws.Cells.Item($r, 4).Value2 = $cVal
}
